# Insert two new weekly data rows ("Primera"/"Segunda" quality pair) for
# Betarraga just before the existing row for 2023-07-21, shifting the
# remainder of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 981:982 (existing rows 981..1053 shift to 983..1055).
$ws.Rows("981:982").Insert()

# New row 981 - "Primera" quality, date 2023-12-05 (serial 45265)
$ws.Cells.Item(981,1).Value = 9
$ws.Cells.Item(981,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(981,3).Value = 'Metropolitana'
$ws.Cells.Item(981,4).Value = 45265
$ws.Cells.Item(981,5).Value = 13
$ws.Cells.Item(981,6).Value = 100114014
$ws.Cells.Item(981,7).Value = 'Betarraga'
$ws.Cells.Item(981,8).Value = 'Sin especificar'
$ws.Cells.Item(981,9).Value = 'Primera'
$ws.Cells.Item(981,10).Value = 5200
$ws.Cells.Item(981,11).Value = 90
$ws.Cells.Item(981,12).Value = 100
$ws.Cells.Item(981,13).Value = 95
$ws.Cells.Item(981,14).Value = '$/unidad'
$ws.Cells.Item(981,15).Value = 'Región Metropolitana'
$ws.Cells.Item(981,16).Value = 95
$ws.Cells.Item(981,17).Value = 1
$ws.Cells.Item(981,18).Value = 'Hortaliza'

# New row 982 - "Segunda" quality, date 2023-12-05 (serial 45265)
$ws.Cells.Item(982,1).Value = 9
$ws.Cells.Item(982,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(982,3).Value = 'Metropolitana'
$ws.Cells.Item(982,4).Value = 45265
$ws.Cells.Item(982,5).Value = 13
$ws.Cells.Item(982,6).Value = 100114014
$ws.Cells.Item(982,7).Value = 'Betarraga'
$ws.Cells.Item(982,8).Value = 'Sin especificar'
$ws.Cells.Item(982,9).Value = 'Segunda'
$ws.Cells.Item(982,10).Value = 3400
$ws.Cells.Item(982,11).Value = 80
$ws.Cells.Item(982,12).Value = 80
$ws.Cells.Item(982,13).Value = 80
$ws.Cells.Item(982,14).Value = '$/unidad'
$ws.Cells.Item(982,15).Value = 'Región Metropolitana'
$ws.Cells.Item(982,16).Value = 80
$ws.Cells.Item(982,17).Value = 1
$ws.Cells.Item(982,18).Value = 'Hortaliza'

# Ensure the date cells keep the date number format used by the rest of
# column D (style carried over from Insert, but set explicitly to be safe).
$ws.Range("D981:D982").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Inserted rows 981-982; sheet now spans $($ws.UsedRange.Address())"
